{"js": "// Rename the \"Referenser\"/\"Kontakt\" label fields from English to Swedish,\n// drop the \", DITE\" suffix from Sven Johansson's position, and update the\n// organisation / work-email / private-email abbreviations throughout the\n// document (both the free-text paragraphs and their mirrored two-column\n// tables).\nconst body = context.document.body;\n\n// The two standalone table-cell values that read simply \"work\" / \"private\"\n// (the email handles in Mattias Schertell's contact table) become \"@work\" /\n// \"@home\". Do this FIRST, while \"work\"/\"private\" are still unambiguous,\n// identifying the right hit by requiring its *entire enclosing paragraph*\n// (i.e. the whole table cell) to equal the search word exactly \u2014 this skips\n// the \"work\" inside \"work email\" and inside \"Work email:work\" without\n// relying on word-boundary matching (which would also treat \"@work\" as\n// containing a whole \"work\" word later on).\nconst wordTargets = [\n  [\"work\", \"@work\"],\n  [\"private\", \"@home\"],\n];\n\nfor (const [searchText, replaceText] of wordTargets) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items/text\");\n  await context.sync();\n  for (const result of results.items) {\n    const paras = result.paragraphs;\n    paras.load(\"items/text\");\n    await context.sync();\n    if (paras.items.length === 1 && paras.items[0].text === searchText) {\n      result.insertText(replaceText, \"Replace\");\n      await context.sync();\n    }\n  }\n}\n\n// Each pair is applied to EVERY matching occurrence in the document (the\n// same literal text appears once in a \"label:value\" paragraph run and once\n// more, split across two table cells, in the table directly below it).\nconst replacements = [\n  // \"Position:Prefekt, DITE\" (paragraph run) and \"Prefekt, DITE\" (table\n  // cell) both lose the \", DITE\" suffix.\n  [\"Prefekt, DITE\", \"Prefekt\"],\n\n  // English -> Swedish field labels (appear for both Sven Johansson and\n  // Markus Fiedler, in both the paragraph run and the table).\n  [\"Organization:\", \"Organisation:\"],\n  [\"Blekinge Institute of Technology, Sweden\", \"DITE, Blekinge Tekniska H\u00f6gskola\"],\n  [\"Email:\", \"Epost:\"],\n  [\"work email\", \"arbete\"],\n  [\"Phone 1:\", \"Tel 1:\"],\n  [\"Phone 2:\", \"Tel 2:\"],\n\n  // Mattias Schertell's contact paragraph: the label+value are glued\n  // together in one run, so replace that combination first (uniquely) ...\n  [\"Work email:work\", \"Epost arbete:@work\"],\n  [\"Private email:private\", \"Epost privat:@home\"],\n  // ... then the corresponding table still has its own separate label and\n  // value cells left to rename.\n  [\"Work email:\", \"Epost arbete:\"],\n  [\"Private email:\", \"Epost privat:\"],\n  [\"Phone:\", \"Tel:\"],\n  [\"Visit Page\", \"bes\u00f6k min profil\"],\n];\n\nfor (const [searchText, replaceText] of replacements) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const result of results.items) {\n    result.insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Rename the \"Referenser\"/\"Kontakt\" label fields from English to Swedish,\n# drop the \", DITE\" suffix from Sven Johansson's position, and update the\n# organisation / work-email / private-email abbreviations throughout the\n# document (both the free-text paragraphs and their mirrored two-column\n# tables).\n#\n# wdFindContinue = 1, wdReplaceAll = 2 (Find.Execute's Wrap / Replace args).\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# The \"Kontakt\" table (3rd table, for Mattias Schertell) has two standalone\n# cell values that read simply \"work\" / \"private\" (the email handles). These\n# are common words that also appear inside \"work email\" (Sven/Markus'\n# tables) and inside \"Work email:work\" below, so rather than a fragile\n# whole-word Find/Replace, set each cell's text directly to \"@work\"/\"@home\".\n$contactTable = $d.Tables.Item(3)\n$contactTable.Cell(2, 2).Range.Text = \"@work\"\n$contactTable.Cell(3, 2).Range.Text = \"@home\"\n\n# Each pair below is applied to EVERY matching occurrence in the document\n# (the same literal text appears once in a \"label:value\" paragraph run and\n# once more, split across two table cells, in the table directly below it).\n\n# \"Position:Prefekt, DITE\" (paragraph run) and \"Prefekt, DITE\" (table cell)\n# both lose the \", DITE\" suffix.\nReplace-AllText \"Prefekt, DITE\" \"Prefekt\"\n\n# English -> Swedish field labels (appear for both Sven Johansson and\n# Markus Fiedler, in both the paragraph run and the table).\nReplace-AllText \"Organization:\" \"Organisation:\"\nReplace-AllText \"Blekinge Institute of Technology, Sweden\" \"DITE, Blekinge Tekniska H\u00f6gskola\"\nReplace-AllText \"Email:\" \"Epost:\"\nReplace-AllText \"work email\" \"arbete\"\nReplace-AllText \"Phone 1:\" \"Tel 1:\"\nReplace-AllText \"Phone 2:\" \"Tel 2:\"\n\n# Mattias Schertell's contact paragraph: the label+value are glued together\n# in one run, so replace that combination (uniquely) in one shot \u2014 the\n# table's label cell (now just \"Work email:\"/\"Private email:\", since its\n# value cell was already renamed above) is handled right after.\nReplace-AllText \"Work email:work\" \"Epost arbete:@work\"\nReplace-AllText \"Private email:private\" \"Epost privat:@home\"\nReplace-AllText \"Work email:\" \"Epost arbete:\"\nReplace-AllText \"Private email:\" \"Epost privat:\"\nReplace-AllText \"Phone:\" \"Tel:\"\nReplace-AllText \"Visit Page\" \"bes\u00f6k min profil\"\n"}
